$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual relay setting values (rows 2-4)
$ws.Range("H2").Value = 0
$ws.Range("G3").Value = 600
$ws.Range("H3").Value = 0
$ws.Range("G4").Value = 600

# Remove the data rows 5-10 (their content becomes obsolete for this
# debugging pass); deleting them shifts the existing blank template rows
# (previously 11-19) up into position, matching the trimmed table.
$ws.Rows("5:10").Delete()

# Move the notes/license textbox up: it keeps the same offset within its
# anchor row (1.5pt / 19050 EMU), but the anchor row moves from 16 to 10
# (0-based), i.e. 1-based row 17 -> row 11. Keep the same sub-row offset.
$shp = $ws.Shapes.Item(1)
$rowOffPt = 19050 / 12700
$shp.Top = $ws.Rows("1:10").Height + $rowOffPt

# Update the active selection to reflect where editing focus ended up
$ws.Range("I6").Select()
